# Add a "Save" column (H) to the s_vals sheet.
# Header H1 mirrors the existing header formatting (bold/bordered/
# centered) by copying G1's format, and H2:H21 get the per-row save
# counts (all 0 except row 12, which is a save = 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy G1's format onto H1 so the new header matches the other headers
# (TB/d2S/K/IP/Win/sum), then set its text.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Save values for rows 2-21 (0 everywhere except row 12, which is 1)
$saveValues = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
